$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Imputed")

# New row 23 of data (ID=22, Type=Staff, ...)
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Staff"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 16.5
$ws.Range("E23").Value = 13.3
$ws.Range("F23").Value = 13.3
$ws.Range("G23").Value = 16.2
$ws.Range("H23").Value = 16.7
$ws.Range("I23").Value = 0.3
$ws.Range("J23").Value = 0.3
$ws.Range("K23").Value = 16.6
$ws.Range("L23").Value = 0.5
$ws.Range("M23").Value = 0.3
$ws.Range("N23").Value = 0.3
$ws.Range("O23").Value = 0.2
$ws.Range("P23").Value = 0.2
$ws.Range("Q23").Value = 7.5
$ws.Range("R23").Value = 3.7
$ws.Range("S23").Value = 0.2
$ws.Range("T23").Value = 12.9
$ws.Range("U23").Value = 0
$ws.Range("W23").Value = 4

# Move selection / top-left view to match the saved worksheet view state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("X23").Select()
